# Update cryptocurrency price/volume snapshot (GitHub Actions scheduled refresh).
# Source data is stored as plain text (prices/percentages keep their original
# formatting, e.g. trailing zeros and thousand-separator dots), so every write
# below is prefixed with a literal apostrophe to force Excel's "text" quote-
# prefix semantics and stop it from re-parsing the string as a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '''64.449.31'
$ws.Range("E2").Value = '''  -2.38%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '''3.180.73'
$ws.Range("E3").Value = '''  -3.93%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '''  -0.01%  '

# Row 5: BNB
$ws.Range("D5").Value = '''571.44'
$ws.Range("E5").Value = '''  -2.39%  '

# Row 6: Solana
$ws.Range("D6").Value = '''169.05'
$ws.Range("E6").Value = '''  -7.33%  '

# Row 8: USDC
$ws.Range("E8").Value = '''  -0.15%  '

# Row 9: LidoStakedEther
$ws.Range("D9").Value = '''3.189.50'
$ws.Range("E9").Value = '''  -3.62%  '

# Row 10: Dogecoin
$ws.Range("E10").Value = '''  -3.63%  '

# Row 11: Toncoin
$ws.Range("D11").Value = '''6.80'
$ws.Range("E11").Value = '''  -0.30%  '

# Row 12: Cardano
$ws.Range("D12").Value = '''0.389'
$ws.Range("E12").Value = '''  -2.85%  '

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = '''3.738.81'
$ws.Range("E13").Value = '''  -3.83%  '

# Row 14: TRON
$ws.Range("E14").Value = '''  -1.60%  '

# Row 15: WrappedBTC
$ws.Range("D15").Value = '''64.501.35'

# Row 16: Avalanche
$ws.Range("D16").Value = '''25.34'
$ws.Range("E16").Value = '''  -2.97%  '

# Row 17: ShibaInu
$ws.Range("E17").Value = '''  -3.45%  '

# Row 18: WrappedEther
$ws.Range("D18").Value = '''3.176.59'
$ws.Range("E18").Value = '''  -4.04%  '

# Row 19: BitcoinCash
$ws.Range("D19").Value = '''419.51'
$ws.Range("E19").Value = '''  -1.12%  '

# Row 20: Chainlink
$ws.Range("D20").Value = '''12.97'
$ws.Range("E20").Value = '''  -1.16%  '

# Row 21: Polkadot
$ws.Range("D21").Value = '''5.37'
$ws.Range("E21").Value = '''  -3.08%  '

# Row 22: Uniswap
$ws.Range("D22").Value = '''7.13'
$ws.Range("E22").Value = '''  -3.36%  '

# Row 23: Dai
$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '''  -0.05%  '

# Row 24: Litecoin
$ws.Range("D24").Value = '''70.38'
$ws.Range("E24").Value = '''  -1.88%  '

# Row 25: LEO
$ws.Range("D25").Value = '''5.68'
$ws.Range("E25").Value = '''  +0.07%  '

# Row 26: Kaspa
$ws.Range("E26").Value = '''  +2.63%  '

# Row 27: Polygon
$ws.Range("D27").Value = '''0.490'
$ws.Range("E27").Value = '''  -4.10%  '

# Row 28: PEPE
$ws.Range("E28").Value = '''  -6.18%  '

# Row 29: InternetComputer(DFINITY)
$ws.Range("E29").Value = '''  -1.29%  '

# Row 30: Binance-PegBSC-USD
$ws.Range("D30").Value = '''0.996'
$ws.Range("E30").Value = '''  -0.31%  '

# Row 31: PancakeSwap
$ws.Range("E31").Value = '''  -3.33%  '

# Row 32: EthereumClassic
$ws.Range("D32").Value = '''21.77'
$ws.Range("E32").Value = '''  -2.51%  '

# Row 33: USDe
$ws.Range("E33").Value = '''  -0.09%  '

# Row 34: NEARProtocol
$ws.Range("E34").Value = '''  -1.85%  '

# Row 35: Aptos
$ws.Range("D35").Value = '''6.36'
$ws.Range("E35").Value = '''  -2.79%  '

# Row 36: Fetch.AI
$ws.Range("E36").Value = '''  -3.41%  '

# Row 37: Monero
$ws.Range("D37").Value = '''156.49'
$ws.Range("E37").Value = '''  -2.41%  '

# Row 38: ImmutableX
$ws.Range("E38").Value = '''  -4.51%  '

# Row 39: Stacks
$ws.Range("E39").Value = '''  -4.62%  '

# Row 40: Maker
$ws.Range("D40").Value = '''2.695.69'
$ws.Range("E40").Value = '''  -6.10%  '

# Row 41: Filecoin
$ws.Range("E41").Value = '''  -1.76%  '

# Row 42: EnergySwap
$ws.Range("D42").Value = '''24.24'
$ws.Range("E42").Value = '''  -7.63%  '

# Row 43: OKB
$ws.Range("D43").Value = '''39.26'
$ws.Range("E43").Value = '''  -1.47%  '

# Row 44: Mantle
$ws.Range("D44").Value = '''0.716'

# Row 45: Hedera
$ws.Range("E45").Value = '''  -5.49%  '

# Row 46: RenderToken
$ws.Range("D46").Value = '''5.60'
$ws.Range("E46").Value = '''  -4.97%  '

# Row 47: VeChain
$ws.Range("E47").Value = '''  -2.39%  '

# Row 48: Bittensor
$ws.Range("D48").Value = '''291.95'

# Row 49: InjectiveProtocol
$ws.Range("D49").Value = '''21.45'
$ws.Range("E49").Value = '''  -6.84%  '

# Row 50: Stellar -> FirstDigitalUSD (rows 50/51 swap rank order)
$ws.Range("B50").Value = '''FirstDigitalUSD'
$ws.Range("C50").Value = '''https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D50").Value = '''0.998'
$ws.Range("E50").Value = '''  -0.26%  '

# Row 51: FirstDigitalUSD -> Stellar
$ws.Range("B51").Value = '''Stellar'
$ws.Range("C51").Value = '''https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").Value = '''0.0991'
$ws.Range("E51").Value = '''  -5.41%  '
